$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add P1=14, Q1=15, with the same formatting as the other header cells (e.g. O1)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$src = $ws.Range("O1")
$dst = $ws.Range("P1:Q1")
$dst.Borders.LineStyle = $src.Borders.Item(7).LineStyle
$dst.Font.Bold = $src.Font.Bold
$dst.HorizontalAlignment = $src.HorizontalAlignment
$dst.VerticalAlignment = $src.VerticalAlignment

# Data rows 2-25: swap I<->K and M<->O, and append P=2, Q=2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: was 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: was 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: was 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: was 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new column
    $ws.Cells.Item($r, 17).Value = 2  # Q: new column
}
